# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1310
$ws1.Range("F4").Value = 1178
$ws1.Range("F5").Value = 14401
$ws1.Range("F6").Value = 16807
$ws1.Range("F9").Value = 31
$ws1.Range("F11").Value = 204
$ws1.Range("F18").Value = 115
$ws1.Range("F20").Value = 1281
$ws1.Range("F23").Value = 49
$ws1.Range("F25").Value = 3
$ws1.Range("F26").Value = 6881
$ws1.Range("F29").Value = 1143
$ws1.Range("F30").Value = 15
$ws1.Range("F32").Value = 5783
$ws1.Range("F36").Value = 4903

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1310
$ws4.Range("F4").Value = 1178
$ws4.Range("F5").Value = 14401
$ws4.Range("F6").Value = 16807
$ws4.Range("F9").Value = 31
$ws4.Range("F11").Value = 204
$ws4.Range("F18").Value = 115
$ws4.Range("F20").Value = 1281
$ws4.Range("F24").Value = 49
$ws4.Range("F26").Value = 3
$ws4.Range("F27").Value = 6881
$ws4.Range("F30").Value = 1143
$ws4.Range("F31").Value = 15
$ws4.Range("F34").Value = 5783
$ws4.Range("F38").Value = 4903

$wb.Save()
